$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels and Devices")

# --- Step 1: Panel Accessories Devices / Label: IOB800 removed from the
#     combos. Done first so that the snapshot sheet (Sheet2, added below)
#     already reflects this update, matching the saved file. ---
$ws.Range("K8").Value = "FB800,POS800-S,POS800-M"
$ws.Range("L8").Value = "Fuse board - 1,POS800-S - 1,POS800-M - 1"

# --- Step 2: duplicate "Add Panels and Devices" to the end of the workbook
#     as "Sheet2" -- a point-in-time snapshot (keeps the old User Story and
#     old Expected 3rd/4th 24V PSU Load values, but the new K8/L8 above). ---
$ws.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Sheet2"

# Deselect the freshly-created sheet (select the whole sheet, matching a
# "click away" state) so it isn't left as the active tab.
$newSheet.Cells.Select() | Out-Null

# --- Step 3: remaining data edits on the original sheet. ---

# Expected 3rd/4th 24V PSU Load values updated to reflect the removal.
$ws.Range("N8").Value = 0.815
$ws.Range("O8").Value = 0.815

# User Story updated, and its highlighted formatting is cleared.
$ws.Range("B4").Value = "NGC-571/T1397 OR TC-183"
$ws.Range("B4").ClearFormats()

# --- Step 4: restore "Add Panels and Devices" as the active sheet/tab and
#     position its selection like the saved file. ---
$ws.Select() | Out-Null
$ws.Range("O8").Select() | Out-Null
